$wb = $excel.ActiveWorkbook

# Add a descriptive title in cell A1 (and other blank "section header" rows)
# of each worksheet, matching the sheet's purpose.

$ws1 = $wb.Worksheets.Item("Overall School Data")
$ws1.Range("A1").Value = "Overall School Data"

$ws2 = $wb.Worksheets.Item("District Summary")
$ws2.Range("A1").Value = "District Summary"

$ws3 = $wb.Worksheets.Item("School Performance")
$ws3.Range("A1").Value = "Highest-Performing Schools (by % Overall Passing)"
$ws3.Range("A9").Value = "Lowest-Performing Schools (by % Overall Passing)"

$ws4 = $wb.Worksheets.Item("Scores by Grade")
$ws4.Range("A1").Value = "Math Scores by Grade"
$ws4.Range("A19").Value = "Reading Scores by Grade"

$ws5 = $wb.Worksheets.Item("Scores by School Factors")
$ws5.Range("A1").Value = "Scores by School Spending"
$ws5.Range("A8").Value = "Scores by School Size"
$ws5.Range("A14").Value = "Scores by School Type"
